$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 (player reset to zero stats, name cleared)
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = ""
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0

# Row 32 (player stats replaced with new summoner)
$ws.Range("B32").Value = 1.297897196261682
$ws.Range("C32").Value = 2222
$ws.Range("D32").Value = 0.002628504672897196
$ws.Range("E32").Value = 4.5
$ws.Range("F32").Value = 272.5
$ws.Range("G32").Value = "Cevahir Akkanat"
$ws.Range("I32").Value = 0.1591705607476636
$ws.Range("J32").Value = 23.5
$ws.Range("K32").Value = 0.01372663551401869

# Row 33 (player reset to zero stats, name cleared, role set to NONE)
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = ""
$ws.Range("H33").Value = "NONE"
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0

# Row 34 (id renumbered, stats reset to zero, name cleared)
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = ""
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0

# Row 35 removed entirely
$ws.Rows("35").Delete()
